# Update the "合肥-漫展信息" workbook: the first 5 upcoming-event rows (old rows 2-6)
# have passed / been superseded, so they are removed and the remaining events shift
# up. A few "want to go" counts (column F) were also refreshed. The same edit is
# applied to both the "展览" sheet and the duplicate "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# F-column ("想去人数") refreshes, keyed by the *new* row number (after the old
# rows 2-6 have been removed and everything below shifts up by 5).
$fUpdates = @{
    4  = 6148
    5  = 167
    6  = 20
    8  = 1853
    9  = 1385
    11 = 931
    12 = 190
    13 = 5555
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Sheets.Item($sheetName)

    # Remove the five events that already happened / were replaced.
    $ws.Rows("2:6").Delete()

    # Column A is a plain sequential index (1..13), independent of which event
    # ended up in that row, so re-stamp it after the shift.
    for ($i = 2; $i -le 14; $i++) {
        $ws.Cells.Item($i, 1).Value = $i - 1
    }

    # Refresh the "want to go" counts that changed.
    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }
}
